$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 currently holds the "old" record (20-10-2021, vol 600, ...).
# The edit: row 16 gets new values (2022-11-11 record), and a brand-new
# row 17 is appended holding the previous row-16 values unchanged.

# 1) Copy the existing row 16 values down into the new row 17 first,
#    so nothing is lost when row 16 is overwritten.
$ws.Range("A17").Value = $ws.Range("A16").Value2
$ws.Range("B17").Value = $ws.Range("B16").Value2
$ws.Range("C17").Value = $ws.Range("C16").Value2
$ws.Range("D17").Value = $ws.Range("D16").Value2
$ws.Range("D17").NumberFormat = $ws.Range("D16").NumberFormat
$ws.Range("E17").Value = $ws.Range("E16").Value2
$ws.Range("F17").Value = $ws.Range("F16").Value2
$ws.Range("G17").Value = $ws.Range("G16").Value2
$ws.Range("H17").Value = $ws.Range("H16").Value2
$ws.Range("I17").Value = $ws.Range("I16").Value2
$ws.Range("J17").Value = $ws.Range("J16").Value2
$ws.Range("K17").Value = $ws.Range("K16").Value2
$ws.Range("L17").Value = $ws.Range("L16").Value2
$ws.Range("M17").Value = $ws.Range("M16").Value2
$ws.Range("N17").Value = $ws.Range("N16").Value2
$ws.Range("O17").Value = $ws.Range("O16").Value2
$ws.Range("P17").Value = $ws.Range("P16").Value2
$ws.Range("Q17").Value = $ws.Range("Q16").Value2
$ws.Range("R17").Value = $ws.Range("R16").Value2

# 2) Overwrite row 16 with the new record values.
$ws.Range("D16").Value = 44876
$ws.Range("J16").Value = 350
$ws.Range("K16").Value = 1500
$ws.Range("L16").Value = 1600
$ws.Range("M16").Value = 1557
$ws.Range("P16").Value = 1557
